$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of rows whose B:AC content must be swapped with one another
# (column A, the sequential row index, stays untouched).
$pairs = @(
    @(29, 30),
    @(36, 37),
    @(49, 50),
    @(76, 77),
    @(87, 88),
    @(111, 112),
    @(122, 123)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1" + ":AC$r1")
    $range2 = $ws.Range("B$r2" + ":AC$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
